$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 186.6960269479135
$ws.Range("I2").Value = 31.60971188545227
$ws.Range("K2").Value = 37
$ws.Range("L2").Value = 0.6250641345977783
$ws.Range("M2").Value = 50.57034972226132
$ws.Range("P2").Value = 0.994824331729891
$ws.Range("R2").Value = 0.9725808417383324
$ws.Range("T2").Value = 0.994824331729891
$ws.Range("V2").Value = 0.9658255928677748
$ws.Range("X2").Value = 0.9795594864329188
$ws.Range("H3").Value = 537.8993470910575
$ws.Range("I3").Value = 35.9458270072937
$ws.Range("K3").Value = 13
$ws.Range("L3").Value = 0.6741127967834473
$ws.Range("M3").Value = 53.32316368834781
$ws.Range("N3").Value = 0.001579694065001903
$ws.Range("O3").Value = 0.02562779079587741
$ws.Range("P3").Value = 0.9938500882908117
$ws.Range("Q3").Value = 0.9943981002252938
$ws.Range("R3").Value = 0.9667621574284866
$ws.Range("S3").Value = 0.969813609955933
$ws.Range("T3").Value = 0.9938500882908117
$ws.Range("U3").Value = 0.9943981002252938
$ws.Range("V3").Value = 0.9692731960366257
$ws.Range("W3").Value = 0.97092948366718
$ws.Range("X3").Value = 0.9642809534667531
$ws.Range("Y3").Value = 0.9687036289968954
$ws.Range("H4").Value = 1594.280011506716
$ws.Range("I4").Value = 69.56055402755737
$ws.Range("L4").Value = 1.209961891174316
$ws.Range("M4").Value = 57.48987181740581
$ws.Range("H5").Value = 49.41097724230254
$ws.Range("I5").Value = 0.4972898960113525
$ws.Range("L5").Value = 0.2443199157714844
$ws.Range("M5").Value = 2.03540466376255
$ws.Range("H6").Value = 138.1121856866538
$ws.Range("I6").Value = 0.9521069526672363
$ws.Range("L6").Value = 0.08733320236206055
$ws.Range("M6").Value = 10.90200435706057
$ws.Range("N6").Value = 0.001431014891680268
$ws.Range("O6").Value = 0.02811760991289549
$ws.Range("Q6").Value = 0.8285714285714286
$ws.Range("S6").Value = 0.5915966386554622
$ws.Range("U6").Value = 0.8285714285714286
$ws.Range("W6").Value = 0.5634920634920635
$ws.Range("Y6").Value = 0.6288561936402468
$ws.Range("H7").Value = 406.3967264224473
$ws.Range("I7").Value = 5.503466844558716
$ws.Range("L7").Value = 0.2481122016906738
$ws.Range("M7").Value = 22.18136313755335
$ws.Range("H13").Value = 155.3291394061399
$ws.Range("I13").Value = 2.278520107269287
$ws.Range("L13").Value = 0.2557082176208496
$ws.Range("M13").Value = 8.91062527621914
$ws.Range("H14").Value = 474.8406113537118
$ws.Range("I14").Value = 6.829186916351318
$ws.Range("J14").Value = $false
$ws.Range("K14").Value = 100
$ws.Range("L14").Value = 0.2205660343170166
$ws.Range("M14").Value = 30.96209685003367
$ws.Range("N14").Value = 0.00158824764965587
$ws.Range("O14").Value = 0.02595488759211882
$ws.Range("P14").Value = 0.9583333333333334
$ws.Range("Q14").Value = 0.9411764705882353
$ws.Range("R14").Value = 0.958266450502248
$ws.Range("S14").Value = 0.9410440911633036
$ws.Range("T14").Value = 0.9583333333333334
$ws.Range("U14").Value = 0.9411764705882353
$ws.Range("V14").Value = 0.9582971178249835
$ws.Range("W14").Value = 0.9414496448734293
$ws.Range("X14").Value = 0.9582370684122705
$ws.Range("Y14").Value = 0.9407649852393787
$ws.Range("H15").Value = 753.1784827059049
$ws.Range("I15").Value = 14.05977082252502
$ws.Range("L15").Value = 0.2785389423370361
$ws.Range("M15").Value = 50.47685865595231
